# Generate Report for Handoff
#
# This script replays a "re-run of the localization handoff/report job" on
# localization-status.xlsx: a new source markdown file
# (75a75d84-8f26-439a-8e59-0bd724a0616e.md, replacing
#  4abcd790-a212-4638-83b9-abe5beef925b.md) was handed off again, producing
# fresh xliff packages (new content hash, new timestamps) and clearing out
# the stale "Latest Target File" / "Latest Handback File" / datetime columns
# until a handback actually happens.

$wb = $excel.ActiveWorkbook

$oldId  = "4abcd790-a212-4638-83b9-abe5beef925b"
$newId  = "75a75d84-8f26-439a-8e59-0bd724a0616e"
$newHash = "2d7908a2f98b9a87041938ded6f18679f4312bfa"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-08-25 06:58:37"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-25 06:58:32"

# Target/handback for this generation haven't happened yet -> cleared out.
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Font.Underline = $false
$wsZhCn.Range("I2").Font.ColorIndex = 0

$wsZhCn.Range("J2").Value = ""

$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Columns I/J shrink now that they no longer hold long file names.
$wsZhCn.Columns.Item(9).ColumnWidth = 17.8
$wsZhCn.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-25 06:58:37"

$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Font.Underline = $false
$wsDeDe.Range("I2").Font.ColorIndex = 0

$wsDeDe.Range("J2").Value = ""

$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Columns.Item(9).ColumnWidth = 17.8
$wsDeDe.Columns.Item(10).ColumnWidth = 20.8
